$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-04 06:37:39"

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
